$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values re-pulled/re-pushed from source data and mean calculation
$ws.Range("F2").Value = -2
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = 3
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = 3
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = 2
$ws.Range("F14").Value = -1
$ws.Range("F15").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = 5
